$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values to the new string identifiers (these contain letters,
# so Excel stores them as text without needing a number-format change)
$ws.Range("C2").Value = "WDPE059A64"
$ws.Range("C3").Value = "WLPE058669A"

# Move the active selection to C3, matching the saved selection state
$ws.Range("C3").Select()
